$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (row 2 through 28) holds a "last changed" date serial number.
# Bump it from 45527 (2024-08-23) to 45528 (2024-08-24) for all rows.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45527) {
        $cell.Value = 45528
    }
}
